$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the value shared by A1:A4 with the new text, then remove the old row 5
$ws.Range("A1").Value = "Paris 2023 Challengers Sticker Capsule"
$ws.Range("A2").Value = "Paris 2023 Challengers Sticker Capsule"
$ws.Range("A3").Value = "Paris 2023 Challengers Sticker Capsule"
$ws.Range("A4").Value = "Paris 2023 Challengers Sticker Capsule"

$ws.Range("A5").Delete()
